$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32-46: Detected Category (column C) GEOGRAPHY_QA -> GEOSPATIAL_QA
for ($r = 32; $r -le 46; $r++) {
    $ws.Cells.Item($r, 3).Value = "GEOSPATIAL_QA"
}

# Rows 47-61: Expected Category (column B) VISUAL_QA -> BINARY_VISUAL_QA
for ($r = 47; $r -le 61; $r++) {
    $ws.Cells.Item($r, 2).Value = "BINARY_VISUAL_QA"
}

# Rows 47-61: Detected Category (column C) specific updates
$ws.Cells.Item(48, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(49, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(50, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(52, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(53, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(54, 3).Value = "IMAGE_RETRIEVAL_BY_IMAGE"
$ws.Cells.Item(57, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(58, 3).Value = "IMAGE_RETRIEVAL_BY_IMAGE"
$ws.Cells.Item(59, 3).Value = "IMAGE_RETRIEVAL_BY_IMAGE"
$ws.Cells.Item(60, 3).Value = "BINARY_VISUAL_QA"
$ws.Cells.Item(61, 3).Value = "BINARY_VISUAL_QA"
